$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Stat Drill" row (day 14 / row 23) complete, like the rows above it.
$ws.Range("C23:J23").Value = "x"

# New "simple past" verb drill table in columns Q:S (rows 24-42).
# Column Q: combined "infinitive (meaning) - simple past" strings, entered
# in the order they were originally typed.
$ws.Range("Q24").Value = "sein (to be) - war"
$ws.Range("Q29").Value = "haben (to have) - hatte"
$ws.Range("Q40").Value = "werden (to become) - wurde"
$ws.Range("Q31").Value = "können (can) - konnte"
$ws.Range("Q35").Value = "müssen (must) - musste"
$ws.Range("Q36").Value = "sagen (to say) - sagte"
$ws.Range("Q34").Value = "machen (to make/do) - machte"
$ws.Range("Q27").Value = "geben (to give) - gab"
$ws.Range("Q30").Value = "kommen (to come) - kam"
$ws.Range("Q38").Value = "sollen (should) - sollte"
$ws.Range("Q42").Value = "wollen (to want) - wollte"
$ws.Range("Q28").Value = "gehen (to go) - ging"
$ws.Range("Q41").Value = "wissen (to know) - wusste"
$ws.Range("Q37").Value = "sehen (to see) - sah"
$ws.Range("Q32").Value = "lassen (to let/allow) - ließ"
$ws.Range("Q39").Value = "stehen (to stand) - stand"
$ws.Range("Q26").Value = "finden (to find) - fand"
$ws.Range("Q25").Value = "bleiben (to stay) - blieb"
$ws.Range("Q33").Value = "liegen (to lie) - lag"

# Columns R (infinitive) and S (simple past form), filled row by row in the
# final, alphabetised order.
$ws.Range("R24").Value = "sein (to be)"
$ws.Range("S24").Value = "war"
$ws.Range("R25").Value = "bleiben (to stay)"
$ws.Range("S25").Value = "blieb"
$ws.Range("R26").Value = "finden (to find)"
$ws.Range("S26").Value = "fand"
$ws.Range("R27").Value = "geben (to give)"
$ws.Range("S27").Value = "gab"
$ws.Range("R28").Value = "gehen (to go)"
$ws.Range("S28").Value = "ging"
$ws.Range("R29").Value = "haben (to have)"
$ws.Range("S29").Value = "hatte"
$ws.Range("R30").Value = "kommen (to come)"
$ws.Range("S30").Value = "kam"
$ws.Range("R31").Value = "können (can)"
$ws.Range("S31").Value = "konnte"
$ws.Range("R32").Value = "lassen (to let/allow)"
$ws.Range("S32").Value = "ließ"
$ws.Range("R33").Value = "liegen (to lie)"
$ws.Range("S33").Value = "lag"
$ws.Range("R34").Value = "machen (to make/do)"
$ws.Range("S34").Value = "machte"
$ws.Range("R35").Value = "müssen (must)"
$ws.Range("S35").Value = "musste"
$ws.Range("R36").Value = "sagen (to say)"
$ws.Range("S36").Value = "sagte"
$ws.Range("R37").Value = "sehen (to see)"
$ws.Range("S37").Value = "sah"
$ws.Range("R38").Value = "sollen (should)"
$ws.Range("S38").Value = "sollte"
$ws.Range("R39").Value = "stehen (to stand)"
$ws.Range("S39").Value = "stand"
$ws.Range("R40").Value = "werden (to become)"
$ws.Range("S40").Value = "wurde"
$ws.Range("R41").Value = "wissen (to know)"
$ws.Range("S41").Value = "wusste"
$ws.Range("R42").Value = "wollen (to want)"
$ws.Range("S42").Value = "wollte"

# Filter applied over the new table (creates the hidden _FilterDatabase name).
$name = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("Q24:Q60"))
$name.Visible = $false

# Selection left on the newly-added table.
$ws.Range("R24:S42").Select()
